# 5.5.1 workbook update:
#  - Rename indicator label in row 5 (national parliaments) and fix the 2010 value
#  - Add a 2022 column (S) with new data for both indicator rows
#  - Add a new row 6 "Proportion of seats held by women in local government*" with its data
#  - Add a new footnote row 7
#  - Update sheet view (drop frozen/topLeft scroll position, change the active selection)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Row 5 fixes: label text, corrected 2010 figure, new 2022 figure
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "Улуттук парламенттердеги аялдардын орундарынын үлүшү"
$ws.Range("G5").Value = 23.9

# New column S, year 2022, modelled on column R
$ws.Range("R4").Copy()
$ws.Range("S4").PasteSpecial(-4122)
$ws.Range("S4").Value = 2022

$ws.Range("R5").Copy()
$ws.Range("S5").PasteSpecial(-4122)
$ws.Range("S5").Value = 21.1

# Row 5 no longer carries the thick bottom rule (it moves down to row 6), but
# it keeps its thick top rule
$ws.Range("A5:S5").Borders.Item(9).LineStyle = -4142
$ws.Range("A5:S5").RowHeight = 27

# ---------------------------------------------------------------------------
# 2. New row 6: local-government indicator
# ---------------------------------------------------------------------------
$ws.Range("A5:S5").Copy()
$ws.Range("A6:S6").PasteSpecial(-4122)

$ws.Range("A6").Value = "Жергиликтүү өз алдынча башкаруу органдарындагы аялдардын орундарынын үлүшү*"
$ws.Range("B6").Value = "Доля мест, занимаемых женщинами в местных органах власти*"
$ws.Range("C6").Value = "Proportion of seats held by women in local government*"

$ws.Range("D6:L6").ClearContents()

$ws.Range("M6").Value = 15.61
$ws.Range("N6").Value = 15.09
$ws.Range("O6").Value = 14.96
$ws.Range("P6").Value = 15.16
$ws.Range("Q6").Value = 14.98
$ws.Range("R6").Value = 31.55
$ws.Range("S6").Value = 36.46
$ws.Range("M6:S6").NumberFormat = "0.0"

$ws.Range("A6:S6").Borders.Item(8).LineStyle = -4142
$ws.Range("A6:S6").Borders.Item(9).Weight = -4138
$ws.Range("A6:S6").RowHeight = 28.5

# ---------------------------------------------------------------------------
# 3. New row 7: footnote
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "*КР ШРӨБК маалыматтары боюнча"
$ws.Range("B7").Value = "*по данным ЦКВПР КР"
$ws.Range("C7").Value = "*according to the CCER of KR"

$ws.Range("A7:C7").Font.Name = "Times New Roman"
$ws.Range("A7:C7").Font.Size = 11
$ws.Range("A7:C7").Font.Bold = $false
$ws.Range("A7:C7").WrapText = $false
$ws.Range("A7:C7").VerticalAlignment = -4160
$ws.Range("A7:C7").Borders.Item(8).LineStyle = -4142
$ws.Range("A7:C7").Borders.Item(9).LineStyle = -4142

# ---------------------------------------------------------------------------
# 4. Sheet view: clear the frozen/left-scrolled column and move the selection
# ---------------------------------------------------------------------------
$ws.Range("T4").Select()
